$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.556.04"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "3.763.31"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D7").Value = "3.762.30"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.17"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "4.392.88"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "3.764.48"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "68.518.10"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  -4.09%  "
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.37"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.07"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("D30").Value = "3.909.80"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.34"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.02"
$ws.Range("D33").ClearFormats()
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").Value = "3.718.99"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.40"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -9.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.303"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.08%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +8.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.73"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.66"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "389.44"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.27%  "
